$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Log_Muestras")

$ws.Range("Z2").Value = "2025-10-17T07:09:30.376317"
$ws.Range("Z3:Z6").Value = "2025-10-17T07:09:30.377314"
$ws.Range("Z7:Z15").Value = "2025-10-17T07:09:30.378314"
$ws.Range("Z16:Z25").Value = "2025-10-17T07:09:30.432590"
$ws.Range("Z26:Z34").Value = "2025-10-17T07:09:30.514926"
$ws.Range("Z35:Z43").Value = "2025-10-17T07:09:30.515925"
$ws.Range("Z44:Z48").Value = "2025-10-17T07:09:30.516929"
